$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 466.33334
$ws.Range("I5").Value = 150
$ws.Range("K5").Value = 150
$ws.Range("M5").Value = -35

$ws.Range("H12").Value = 945.3333
$ws.Range("I12").Value = 2598
$ws.Range("K12").Value = 2598
$ws.Range("M12").Value = -2428

$ws.Range("H17").Value = 919.6667
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 919.6667
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2759.0001
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3095.0001

$ws.Range("H33").Value = 1133
$ws.Range("I33").Value = 449.5
$ws.Range("K33").Value = 449.5
$ws.Range("M33").Value = -220.5

$ws.Range("H55").Value = 221.5
$ws.Range("I55").Value = 168.75
$ws.Range("K55").Value = 168.75
$ws.Range("M55").Value = 45.25

$ws.Range("H98").Value = 797.4286
$ws.Range("I98").Value = 797.4286
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 797.4286
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 700.5714
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 797.4286
$ws.Range("I122").Value = 797.4286
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2392.2858
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 57.71420000000035
$ws.Range("N122").ClearContents()

$ws.Range("H125").Value = 3208.1428
$ws.Range("J125").Value = 3383.75
$ws.Range("L125").Value = 30453.75
$ws.Range("N125").Value = -35373.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 250
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -134
$ws.Range("N4").ClearContents()

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H74").Value = 799.4
$ws.Range("I74").Value = 799.3333
$ws.Range("K74").Value = 799.3333
$ws.Range("M74").Value = 74.66669999999999

$ws.Range("H77").Value = 799.4
$ws.Range("I77").Value = 799.3333
$ws.Range("K77").Value = 3996.6665
$ws.Range("M77").Value = 371.3334999999997

$ws.Range("H122").Value = 950
$ws.Range("I122").Value = 950
$ws.Range("K122").Value = 2850
$ws.Range("M122").Value = -400

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2104.3684
$ws.Range("I7").Value = 1358.6666
$ws.Range("K7").Value = 1358.6666
$ws.Range("M7").Value = -1245.6666

$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()

$ws.Range("H31").Value = 2304.2
$ws.Range("J31").Value = 2829.8
$ws.Range("L31").Value = 2829.8
$ws.Range("N31").Value = -3419.8

$ws.Range("H34").Value = 2304.2
$ws.Range("J34").Value = 2829.8
$ws.Range("L34").Value = 2829.8
$ws.Range("N34").Value = -3233.8

$ws.Range("H106").Value = 23799.6
$ws.Range("J106").Value = 23799.6
$ws.Range("L106").Value = 23799.6
$ws.Range("N106").Value = -26323.6

$ws.Range("H134").Value = 1529.75
$ws.Range("I134").Value = 1556.8823
$ws.Range("K134").Value = 4670.6469
$ws.Range("M134").Value = -2135.6469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 173.66667
$ws.Range("I33").Value = 158.83333
$ws.Range("J33").Value = 203.33333
$ws.Range("K33").Value = 952.9999799999999
$ws.Range("L33").Value = 1219.99998
$ws.Range("M33").Value = -669.9999799999999
$ws.Range("N33").Value = -1785.99998

$ws.Range("H38").Value = 146.7
$ws.Range("I38").Value = 58.375
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 175.125
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = 171.875
$ws.Range("N38").Value = -2194

$ws.Range("H68").Value = 13000
$ws.Range("J68").Value = 13000
$ws.Range("L68").Value = 39000
$ws.Range("N68").Value = -40622

$ws.Range("H71").Value = 13000
$ws.Range("J71").Value = 13000
$ws.Range("L71").Value = 117000
$ws.Range("N71").Value = -125112

$ws.Range("H129").Value = 1791.125
$ws.Range("I129").Value = 783.125
$ws.Range("J129").Value = 2799.125
$ws.Range("K129").Value = 2349.375
$ws.Range("L129").Value = 8397.375
$ws.Range("M129").Value = 2650.625
$ws.Range("N129").Value = -18397.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 17500002
$ws.Range("I11").Value = 17500002
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 17500002
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -17499863
$ws.Range("N11").ClearContents()

$ws.Range("H20").Value = 7000
$ws.Range("J20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("N20").Value = -7490

$ws.Range("H80").Value = 8253
$ws.Range("J80").Value = 8253
$ws.Range("L80").Value = 8253
$ws.Range("N80").Value = -10249

$ws.Range("H83").Value = 8253
$ws.Range("J83").Value = 8253
$ws.Range("L83").Value = 41265
$ws.Range("N83").Value = -51249

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2200
$ws.Range("I7").Value = 2200
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2200
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2088
$ws.Range("N7").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws.Range("H23").Value = 19000004
$ws.Range("I23").Value = 19000004
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 19000004
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -18999774
$ws.Range("N23").ClearContents()

$ws.Range("H46").Value = 4694.3335
$ws.Range("I46").Value = 1999.6666
$ws.Range("J46").Value = 5233.2666
$ws.Range("K46").Value = 1999.6666
$ws.Range("L46").Value = 5233.2666
$ws.Range("M46").Value = -1811.6666
$ws.Range("N46").Value = -5609.2666

$ws.Range("H55").Value = 439.1
$ws.Range("I55").Value = 194.5
$ws.Range("J55").Value = 500.25
$ws.Range("K55").Value = 194.5
$ws.Range("L55").Value = 500.25
$ws.Range("M55").Value = -21.5
$ws.Range("N55").Value = -846.25

$ws.Range("H56").Value = 40057
$ws.Range("J56").Value = 40057
$ws.Range("L56").Value = 40057
$ws.Range("M56").Value = -41439

$ws.Range("H101").Value = 19330.666
$ws.Range("J101").Value = 19330.666
$ws.Range("L101").Value = 19330.666
$ws.Range("N101").Value = -25820.666

$ws.Range("H108").Value = 100000
$ws.Range("I108").Value = 100000
$ws.Range("K108").Value = 100000
$ws.Range("M108").Value = -96160

$ws.Range("H126").Value = 2200
$ws.Range("I126").Value = 2200
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6600
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4130
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30001
$ws.Range("J64").Value = 30001
$ws.Range("L64").Value = 30001
$ws.Range("N64").Value = -30497

$ws.Range("H67").Value = 30001
$ws.Range("J67").Value = 30001
$ws.Range("L67").Value = 30001
$ws.Range("N67").Value = -31717

$ws.Range("H107").Value = 859
$ws.Range("J107").Value = 1133.6666
$ws.Range("L107").Value = 3400.9998
$ws.Range("N107").Value = -7240.9998

$ws.Range("H136").Value = 908.82355
$ws.Range("I136").Value = 908.82355
$ws.Range("K136").Value = 2726.47065
$ws.Range("M136").Value = -176.4706499999998
